$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 154 (shifts existing rows 154-172 down to 155-173),
# adding a new weekly price record for Acelga at the top of this block.
$ws.Rows.Item(154).Insert()

$ws.Cells.Item(154, 1).Value = 11
$ws.Cells.Item(154, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(154, 3).Value = "Bíobío"
$ws.Cells.Item(154, 4).Value = 44504
$ws.Cells.Item(154, 5).Value = 8
$ws.Cells.Item(154, 6).Value = 100112009
$ws.Cells.Item(154, 7).Value = "Acelga"
$ws.Cells.Item(154, 8).Value = "Sin especificar"
$ws.Cells.Item(154, 9).Value = "Primera"
$ws.Cells.Item(154, 10).Value = 700
$ws.Cells.Item(154, 11).Value = 600
$ws.Cells.Item(154, 12).Value = 650
$ws.Cells.Item(154, 13).Value = 621
$ws.Cells.Item(154, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(154, 15).Value = "Región de Ñuble"
$ws.Cells.Item(154, 16).Value = 621
$ws.Cells.Item(154, 17).Value = 1
$ws.Cells.Item(154, 18).Value = "Hortaliza"
